$wb = $excel.ActiveWorkbook

# "List" sheet: rename the accountId row/column to owner, and change the
# value expression to use the printer helper.
$listWs = $wb.Worksheets.Item("List")
$listWs.Range("B1").Value = "`${msg.getProperty('document_owner')}"
$listWs.Range("B2").Value = "`${printer.print(document.owner)}"

# "Search" sheet: append a new search-criteria row for "owner".
$searchWs = $wb.Worksheets.Item("Search")
$searchWs.Range("A5").Value = "`${msg.getProperty('document_owner')}"
$searchWs.Range("B5").Value = "`${owner}"
